$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "26.442.48"
Set-TextValue "E2" "  -2.82%  "
Set-TextValue "D3" "1.806.13"
Set-TextValue "E3" "  -2.54%  "
Set-TextValue "E4" "  +0.73%  "
Set-TextValue "E5" "  +0.67%  "
Set-TextValue "D6" "308.04"
Set-TextValue "E6" "  -1.63%  "
Set-TextValue "D7" "0.4557"
Set-TextValue "E7" "  -1.70%  "
Set-TextValue "D8" "0.3661"
Set-TextValue "E8" "  -1.57%  "
Set-TextValue "D9" "0.07126"
Set-TextValue "E9" "  -2.18%  "
Set-TextValue "D10" "0.8763"
Set-TextValue "E10" "  -1.24%  "
Set-TextValue "D11" "0.07788"
Set-TextValue "E11" "  -0.16%  "
Set-TextValue "D12" "19.33"
Set-TextValue "E12" "  -3.57%  "
Set-TextValue "D13" "1.824.02"
Set-TextValue "E13" "  -6.31%  "
Set-TextValue "D14" "5.268"
Set-TextValue "E14" "  -1.96%  "
Set-TextValue "D15" "6.368"
Set-TextValue "E15" "  -2.22%  "
Set-TextValue "D16" "85.83"
Set-TextValue "E16" "  -5.85%  "
Set-TextValue "D17" "1.009"
Set-TextValue "E17" "  +0.78%  "
Set-TextValue "D18" "0.000008570"
Set-TextValue "E18" "  -3.95%  "
Set-TextValue "E19" "  +0.52%  "
Set-TextValue "D20" "26.475.59"
Set-TextValue "E20" "  -2.83%  "
Set-TextValue "D21" "14.24"
Set-TextValue "E21" "  -3.14%  "
Set-TextValue "E22" "  -1.45%  "
Set-TextValue "D23" "10.43"
Set-TextValue "E23" "  -0.60%  "
Set-TextValue "D24" "1.977"
Set-TextValue "E24" "  +1.35%  "
Set-TextValue "D25" "150.91"
Set-TextValue "E25" "  -0.58%  "
Set-TextValue "E26" "  -2.56%  "
Set-TextValue "D27" "2.050"
Set-TextValue "E27" "  +0.38%  "
Set-TextValue "D28" "112.48"
Set-TextValue "E28" "  -2.72%  "
Set-TextValue "D29" "4.836"
Set-TextValue "E29" "  -4.41%  "
Set-TextValue "D30" "0.08665"
Set-TextValue "E30" "  -1.56%  "
Set-TextValue "D31" "3.044"
Set-TextValue "E31" "  -4.24%  "
Set-TextValue "D32" "0.7317"
Set-TextValue "E32" "  -4.38%  "
Set-TextValue "D33" "4.462"
Set-TextValue "E33" "  -0.89%  "
Set-TextValue "E34" "  -5.00%  "
Set-TextValue "D35" "1.005"
Set-TextValue "E35" "  +0.51%  "
Set-TextValue "D36" "2.515"
Set-TextValue "E36" "  -8.27%  "
Set-TextValue "D37" "1.079"
Set-TextValue "E37" "  -0.76%  "
Set-TextValue "D38" "0.01929"
Set-TextValue "E38" "  -0.52%  "
Set-TextValue "D39" "0.05104"
Set-TextValue "E39" "  -2.41%  "
Set-TextValue "D40" "2.890"
Set-TextValue "E40" "  -1.40%  "
Set-TextValue "D41" "6.928"
Set-TextValue "E41" "  -2.00%  "
Set-TextValue "D42" "0.5005"
Set-TextValue "E42" "  -1.93%  "
Set-TextValue "D43" "0.1559"
Set-TextValue "E43" "  -4.19%  "
Set-TextValue "D44" "8.119"
Set-TextValue "E44" "  -3.29%  "
Set-TextValue "D45" "1.008"
Set-TextValue "E45" "  +0.78%  "
Set-TextValue "D46" "0.4607"
Set-TextValue "E46" "  -3.75%  "
Set-TextValue "D47" "9.969"
Set-TextValue "E47" "  -3.50%  "
Set-TextValue "D48" "100.49"
Set-TextValue "E48" "  -2.32%  "
Set-TextValue "D49" "1.585"
Set-TextValue "E49" "  -3.26%  "
Set-TextValue "E50" "  -3.31%  "
Set-TextValue "D51" "63.90"
Set-TextValue "E51" "  -2.27%  "
